$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "231.24", "1.00") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "44.094.81", "  +5.31%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.294.51", "  +2.62%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.07%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "231.24", "  -0.61%  "),
    @(6, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.629", "  +0.06%  "),
    @(7, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "61.01", "  -0.40%  "),
    @(8, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.04%  "),
    @(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.426", "  +4.95%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.0948", "  +3.95%  "),
    @(11, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "57.87", "  -1.90%  "),
    @(12, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.104", "  +0.32%  "),
    @(13, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.636.80", "  +2.68%  "),
    @(14, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "24.25", "  +7.85%  "),
    @(15, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "15.78", "  +0.40%  "),
    @(16, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.91", "  +5.16%  "),
    @(17, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.816", "  +1.28%  "),
    @(18, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.303.90", "  +2.49%  "),
    @(19, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "44.052.91", "  +5.38%  "),
    @(20, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0₃0947", "  +4.07%  "),
    @(21, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "73.74", "  +1.82%  "),
    @(22, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.27", "  +3.62%  "),
    @(23, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "254.43", "  +1.29%  "),
    @(24, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.10%  "),
    @(25, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "2.56", "  +6.79%  "),
    @(26, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.36", "  -0.37%  "),
    @(27, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.93", "  +1.88%  "),
    @(28, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "171.29", "  +1.14%  "),
    @(29, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.140", "  -2.71%  "),
    @(30, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "20.60", "  +2.69%  "),
    @(31, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.42", "  +0.04%  "),
    @(32, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "2.76", "  +1.25%  "),
    @(33, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.122", "  -0.14%  "),
    @(34, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "5.10", "  +0.85%  "),
    @(35, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.80", "  +2.59%  "),
    @(36, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0660", "  +3.31%  "),
    @(37, "THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "6.53", "  -1.81%  "),
    @(38, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "3.65", "  -2.58%  "),
    @(39, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.41", "  +1.80%  "),
    @(40, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0251", "  +4.01%  "),
    @(41, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.00", "  +0.19%  "),
    @(42, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "8.81", "  +2.76%  "),
    @(43, "TerraClassic", "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc", "0.000225", "  -13.29%  "),
    @(44, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "99.30", "  -0.20%  "),
    @(45, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.0969", "  +0.95%  "),
    @(46, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.22", "  -1.51%  "),
    @(47, "Celestia", "https://coinranking.com/coin/YQcD0lBl7+celestia-tia", "10.31", "  +17.13%  "),
    @(48, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "17.05", "  +2.61%  "),
    @(49, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.484.02", "  +0.09%  "),
    @(50, "FTXToken", "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt", "4.39", "  -6.06%  "),
    @(51, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.10", "  +1.38%  ")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
}

# Restore default style so only the cell content changed,
# matching the source diff (no formatting changes).
$ws.Range("D2:D51").Style = "Normal"
